$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new cell values in the exact order that reproduces the original
# shared-string allocation order (first-use order determines the <si> index).
$ws.Range("E11").Value = "Lebanon"
$ws.Range("F11").Value = "SO Tripoli G3"
$ws.Range("A11").Value = "009-80DD3"
$ws.Range("A12").Value = "009-80DD6"
$ws.Range("A13").Value = "009-80DDB"
$ws.Range("F13").Value = "SO Tripoli G2"
$ws.Range("F14").Value = "SO Tripoli G1"
$ws.Range("A14").Value = "009-80DD4"
$ws.Range("A15").Value = "009-80A9E"
$ws.Range("F15").Value = "Lagos"
$ws.Range("H15").Value = "2 generators on 1 GB, not handled yet -- other is 16kVa (not currently used)"

# Remaining cells (reuse existing shared strings / numeric values).
$ws.Range("D11").Value = "Middle East & North Africa"
$ws.Range("D12").Value = "Middle East & North Africa"
$ws.Range("D13").Value = "Middle East & North Africa"
$ws.Range("D14").Value = "Middle East & North Africa"
$ws.Range("D15").Value = "Middle East & North Africa"

$ws.Range("E12").Value = "Lebanon"
$ws.Range("E13").Value = "Lebanon"
$ws.Range("E14").Value = "Lebanon"
$ws.Range("E15").Value = "Nigeria"

$ws.Range("F12").Value = "SO Tripoli G3"

$ws.Range("G11").Value = 200
$ws.Range("G12").Value = 350
$ws.Range("G13").Value = 220
$ws.Range("G14").Value = 80
$ws.Range("G15").Value = 88

# Row 15's meter-serial cell carries a vertically centered style in the
# source workbook (new font + alignment record).
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("A15").Font.Name = "Calibri"

# Match the post-edit selection cursor position (below the last data row).
$ws.Range("H16").Select()
